# Add a new slide (5th) using the "Title and Content" layout (same
# layout used by the deck's other content slides) and fill in the
# Advantages/Disadvantages bullet content.

$p = $ppt.ActivePresentation

$s = $p.Slides.Add(5, 2)

# --- Title -----------------------------------------------------------
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Advantages & Disadvantages"

# --- Body content ------------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Advantages:`rHandle millions of rows`rExcellent speed but low memory consumption`rDisadvantages`rNot perform well on small datasets`rHigh loss on one tree may result in "

# Sub-bullets (second outline level)
$body.Paragraphs(2, 1).IndentLevel = 2
$body.Paragraphs(3, 1).IndentLevel = 2
$body.Paragraphs(5, 1).IndentLevel = 2
$body.Paragraphs(6, 1).IndentLevel = 2

# Last bullet is split across two runs - append the second run so it
# keeps its own run-level formatting, matching the authored deck.
$tail = $body.InsertAfter("less accuracy")
